# Apply updated cryptocurrency price-list values.
# For every touched cell, the source feed writes a value that *looks*
# numeric/percent (e.g. "327.81", "5.79%", "14"), but the sheet stores
# these as plain text (no leading apostrophe shown, no numeric <v> type).
# Forcing NumberFormat="@" (Text) before the write keeps Excel from
# auto-converting them to Number/Percentage cells, then we restore the
# cell style to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "327.81"
Set-TextValue $ws.Range("E2") "5.79%"
Set-TextValue $ws.Range("G2") "14"
# Row 3
Set-TextValue $ws.Range("D3") "39.92"
Set-TextValue $ws.Range("E3") "6.87%"
Set-TextValue $ws.Range("G3") "14"
# Row 4
Set-TextValue $ws.Range("D4") "5.593"
Set-TextValue $ws.Range("E4") "9.00%"
Set-TextValue $ws.Range("G4") "14"
# Row 5
Set-TextValue $ws.Range("D5") "0.08129"
Set-TextValue $ws.Range("E5") "3.69%"
Set-TextValue $ws.Range("G5") "14"
# Row 6
Set-TextValue $ws.Range("D6") "4.555"
Set-TextValue $ws.Range("E6") "3.48%"
Set-TextValue $ws.Range("G6") "14"
# Row 7
Set-TextValue $ws.Range("D7") "8.674"
Set-TextValue $ws.Range("E7") "4.87%"
Set-TextValue $ws.Range("G7") "14"
# Row 8
Set-TextValue $ws.Range("D8") "1.977"
Set-TextValue $ws.Range("E8") "5.23%"
Set-TextValue $ws.Range("G8") "14"
# Row 9
Set-TextValue $ws.Range("G9") "14"
# Row 10
Set-TextValue $ws.Range("D10") "0.9491"
Set-TextValue $ws.Range("E10") "2.65%"
Set-TextValue $ws.Range("G10") "14"
# Row 11
Set-TextValue $ws.Range("D11") "0.1278"
Set-TextValue $ws.Range("E11") "9.50%"
Set-TextValue $ws.Range("G11") "14"
# Row 12
Set-TextValue $ws.Range("D12") "0.1976"
Set-TextValue $ws.Range("E12") "4.19%"
Set-TextValue $ws.Range("G12") "14"
# Row 13
Set-TextValue $ws.Range("D13") "0.09185"
Set-TextValue $ws.Range("E13") "2.88%"
Set-TextValue $ws.Range("G13") "14"
# Row 14
Set-TextValue $ws.Range("D14") "0.03588"
Set-TextValue $ws.Range("E14") "8.03%"
Set-TextValue $ws.Range("G14") "14"
# Row 15
Set-TextValue $ws.Range("D15") "0.09594"
Set-TextValue $ws.Range("E15") "-0.17%"
Set-TextValue $ws.Range("G15") "14"
# Row 16
Set-TextValue $ws.Range("D16") "0.001325"
Set-TextValue $ws.Range("E16") "-3.71%"
Set-TextValue $ws.Range("G16") "14"
# Row 17
Set-TextValue $ws.Range("D17") "0.006159"
Set-TextValue $ws.Range("E17") "-0.68%"
Set-TextValue $ws.Range("G17") "14"
# Row 18
Set-TextValue $ws.Range("D18") "3.371"
Set-TextValue $ws.Range("E18") "-0.63%"
Set-TextValue $ws.Range("G18") "14"
# Row 19
Set-TextValue $ws.Range("D19") "0.3501"
Set-TextValue $ws.Range("E19") "1.27%"
Set-TextValue $ws.Range("G19") "14"
# Row 20
Set-TextValue $ws.Range("D20") "7.454"
Set-TextValue $ws.Range("E20") "16.67%"
Set-TextValue $ws.Range("G20") "14"
# Row 21
Set-TextValue $ws.Range("D21") "0.1361"
Set-TextValue $ws.Range("E21") "3.49%"
Set-TextValue $ws.Range("G21") "14"
# Row 22
Set-TextValue $ws.Range("D22") "0.2488"
Set-TextValue $ws.Range("E22") "3.56%"
Set-TextValue $ws.Range("G22") "14"
# Row 23
Set-TextValue $ws.Range("D23") "0.04426"
Set-TextValue $ws.Range("E23") "1.74%"
Set-TextValue $ws.Range("G23") "14"
# Row 24
Set-TextValue $ws.Range("D24") "0.001225"
Set-TextValue $ws.Range("E24") "2.02%"
Set-TextValue $ws.Range("G24") "14"
# Row 25
Set-TextValue $ws.Range("D25") "0.004287"
Set-TextValue $ws.Range("E25") "0.26%"
Set-TextValue $ws.Range("G25") "14"
# Row 26
Set-TextValue $ws.Range("D26") "0.0001191"
Set-TextValue $ws.Range("E26") "-14.94%"
Set-TextValue $ws.Range("G26") "14"
# Row 27
Set-TextValue $ws.Range("D27") "0.0003993"
Set-TextValue $ws.Range("E27") "37.59%"
Set-TextValue $ws.Range("G27") "14"
# Row 28
Set-TextValue $ws.Range("G28") "14"
# Row 29
Set-TextValue $ws.Range("G29") "14"
# Row 30
Set-TextValue $ws.Range("G30") "14"
# Row 31
Set-TextValue $ws.Range("G31") "14"
# Row 32
Set-TextValue $ws.Range("G32") "14"
# Row 33
Set-TextValue $ws.Range("G33") "14"
# Row 34
Set-TextValue $ws.Range("G34") "14"
# Row 35
Set-TextValue $ws.Range("G35") "14"
# Row 36
Set-TextValue $ws.Range("G36") "14"
# Row 37
Set-TextValue $ws.Range("G37") "14"
# Row 38
Set-TextValue $ws.Range("G38") "14"
# Row 39
Set-TextValue $ws.Range("D39") "0.02525"
Set-TextValue $ws.Range("E39") "16.74%"
Set-TextValue $ws.Range("G39") "14"
# Row 40
Set-TextValue $ws.Range("D40") "0.05212"
Set-TextValue $ws.Range("E40") "3.91%"
Set-TextValue $ws.Range("G40") "14"
# Row 41
Set-TextValue $ws.Range("D41") "0.007741"
Set-TextValue $ws.Range("E41") "2.29%"
Set-TextValue $ws.Range("G41") "14"
# Row 42
Set-TextValue $ws.Range("E42") "5.84%"
Set-TextValue $ws.Range("G42") "14"
# Row 43
Set-TextValue $ws.Range("D43") "0.008847"
Set-TextValue $ws.Range("E43") "4.24%"
Set-TextValue $ws.Range("G43") "14"
# Row 44
Set-TextValue $ws.Range("D44") "0.002192"
Set-TextValue $ws.Range("E44") "5.83%"
Set-TextValue $ws.Range("G44") "14"
# Row 45
Set-TextValue $ws.Range("D45") "0.009621"
Set-TextValue $ws.Range("E45") "18.34%"
Set-TextValue $ws.Range("G45") "14"
# Row 46
Set-TextValue $ws.Range("D46") "0.00006685"
Set-TextValue $ws.Range("E46") "1.63%"
Set-TextValue $ws.Range("G46") "14"
# Row 47
Set-TextValue $ws.Range("E47") "0.02%"
Set-TextValue $ws.Range("G47") "14"
# Row 48
Set-TextValue $ws.Range("D48") "0.002876"
Set-TextValue $ws.Range("E48") "-12.72%"
Set-TextValue $ws.Range("G48") "14"
# Row 49
Set-TextValue $ws.Range("D49") "0.002302"
Set-TextValue $ws.Range("E49") "59.41%"
Set-TextValue $ws.Range("G49") "14"
# Row 50
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "0.02%"
Set-TextValue $ws.Range("G50") "14"
# Row 51
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "0.02%"
Set-TextValue $ws.Range("G51") "14"
